# The Lab Exercise date line ("Lab Exercise 11/30/2023") needs to become
# "Lab Exercise 11/22/2024". The diff shows the change as a series of
# run-local edits ("1/30" -> "1/", "/202" -> "22", "3" -> "/202", plus a new
# trailing run for "4"), but all of those runs share identical bold /
# sz=28 / szCs=28 formatting, so the net visible effect is simply the date
# string changing from 11/30/2023 to 11/22/2024. Do that with a single,
# unambiguous Find/Replace over the whole document.
$d = $word.ActiveDocument

$d.Content.Find.Execute("11/30/2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11/22/2024", 2)
